# Add a new "Wind Onshore" commodity, a new "MIN_EX_WIND_ON" mining process,
# and the corresponding row in the MIN_IMP (mining/import) table.

$wb = $excel.ActiveWorkbook

# --- SEC_Comm: new commodity row (row 9) ---
$wsComm = $wb.Worksheets.Item("SEC_Comm")
$wsComm.Range("B9").Value = "NRG"
$wsComm.Range("C9").Value = "WIND_ON"
$wsComm.Range("D9").Value = "Wind Onshore"
$wsComm.Range("E9").Value = "PJ"
$wsComm.Range("G9").Value = "SEASON"

# --- SEC_Processes: new process row (row 9) ---
$wsProc = $wb.Worksheets.Item("SEC_Processes")
$wsProc.Range("B9").Value = "MIN"
$wsProc.Range("D9").Value = "MIN_EX_WIND_ON"
$wsProc.Range("E9").Value = "Wind mine"
$wsProc.Range("F9").Value = "PJ"
$wsProc.Range("G9").Value = "PJa"
$wsProc.Range("H9").Value = "SEASON"

# --- MIN_IMP: new mining entry row (row 9) referencing the above ---
$wsMinImp = $wb.Worksheets.Item("MIN_IMP")
$wsMinImp.Range("B9").Formula = "=SEC_Processes!D9"
$wsMinImp.Range("C9").Formula = "=SEC_Processes!E9"
$wsMinImp.Range("D9").Formula = "=SEC_Comm!D9"
$wsMinImp.Range("E9").Value = 0.001

# --- Mirror the saved view/selection state on each touched sheet ---
$wsComm.Range("F10").Select() | Out-Null
$wsProc.Range("H10").Select() | Out-Null

$wsMinImp.Activate() | Out-Null
$wsMinImp.Range("E11").Select() | Out-Null
